$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries: Argelia now appears before Noruega ---
# Row 56 held Noruega, row 57 held Argelia; swap the country names so
# row 56 = Argelia, row 57 = Noruega (values for each row updated below).
$ws.Range("A56").Value = "Argelia"
$ws.Range("A57").Value = "Noruega"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 18:35"

# --- Update country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1693157
$ws.Range("C4").Value = 6721
$ws.Range("D4").Value = 453946
$ws.Range("E4").Value = 1139763
$ws.Range("G4").Value = 148
$ws.Range("H4").Value = 99448

# Row 9: Reino Unido
$ws.Range("B9").Value = 230158
$ws.Range("C9").Value = 300
$ws.Range("D9").Value = 141981
$ws.Range("E9").Value = 55300
$ws.Range("G9").Value = 92
$ws.Range("H9").Value = 32877

# Row 11: Francia
$ws.Range("B11").Value = 180583
$ws.Range("C11").Value = 255
$ws.Range("E11").Value = 10992
$ws.Range("G11").Value = 20
$ws.Range("H11").Value = 8391

# Row 16: Canada
$ws.Range("B16").Value = 85104
$ws.Range("C16").Value = 405
$ws.Range("D16").Value = 44207

# Row 53: Chequia
$ws.Range("B53").Value = 9171
$ws.Range("C53").Value = 33
$ws.Range("E53").Value = 4404

# Row 55: Kazajistan
$ws.Range("D55").Value = 4515
$ws.Range("E55").Value = 3981

# Row 56: Argelia (updated figures)
$ws.Range("B56").Value = 8503
$ws.Range("C56").Value = 197
$ws.Range("D56").Value = 4747
$ws.Range("E56").Value = 3147
$ws.Range("G56").Value = 9
$ws.Range("H56").Value = 609

# Row 57: Noruega (unchanged figures, now on row 57)
$ws.Range("B57").Value = 8360
$ws.Range("C57").Value = 8
$ws.Range("D57").Value = 7727
$ws.Range("E57").Value = 398
$ws.Range("H57").Value = 235

# Row 71: Luxemburgo
$ws.Range("D71").Value = 3781
$ws.Range("E71").Value = 102

# Row 101: Eslovaquia
$ws.Range("B101").Value = 1395
$ws.Range("C101").Value = 24
$ws.Range("E101").Value = 1247

# Row 155: Suazilandia
$ws.Range("D155").Value = 71
$ws.Range("E155").Value = 137
